$word.Options.AutoFormatReplaceQuotes = $false
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$d = $word.ActiveDocument
$notFound = @()

$r = $d.Content
$found = $r.Find.Execute('Generated on: 2025-06-16 12:31:55')
if ($found) { $r.Text = 'Generated on: 2025-06-18 13:26:06' } else { $notFound += 'pair_0' }

$r = $d.Content
$found = $r.Find.Execute('**Week 1: Introduction to Python & Basic Syntax**')
if ($found) { $r.Text = '**Week 1: Introduction to Python and Setup**' } else { $notFound += 'pair_1' }

$r = $d.Content
$found = $r.Find.Execute('* **Subtopics:**  Python''s history and uses, installing Python and an IDE (e.g., VS Code, PyCharm), basic program structure (print statements, comments).')
if ($found) { $r.Text = '* **Subtopics:**  Introduction to Python, its applications, and advantages. Installing Python and a suitable IDE (e.g., VS Code, PyCharm).  Basic command-line usage.  Hello, World! program.' } else { $notFound += 'pair_2' }

$r = $d.Content
$found = $r.Find.Execute('* **Activities:**  Install Python and an IDE. Write a "Hello, World!" program.  Practice using comments to document code.')
if ($found) { $r.Text = '* **Activities:**  Install Python and an IDE. Run the "Hello, World!" program. Explore the IDE interface.' } else { $notFound += 'pair_3' }

$r = $d.Content
$found = $r.Find.Execute('**Week 2: Variables, Data Types & Operators**')
if ($found) { $r.Text = '**Week 2: Data Types and Operators**' } else { $notFound += 'pair_4' }

$r = $d.Content
$found = $r.Find.Execute('* **Topic:** Working with Data.')
if ($found) { $r.Text = '* **Topic:** Working with fundamental data types and operators.' } else { $notFound += 'pair_5' }

$r = $d.Content
$found = $r.Find.Execute('* **Subtopics:** Integers, floats, strings, booleans, basic arithmetic operators (+, -, *, /, //, %, **), assignment operators (=, +=, -=, etc.). Type conversion.')
if ($found) { $r.Text = '* **Subtopics:**  Integers, floats, strings, booleans. Arithmetic, comparison, logical, and assignment operators. Type conversion.' } else { $notFound += 'pair_6' }

$r = $d.Content
$found = $r.Find.Execute('* **Activities:**  Write a program to perform calculations.  Convert between different data types.  Practice string manipulation (concatenation, slicing).')
if ($found) { $r.Text = '* **Activities:**  Practice problems involving different data types and operators.  Simple calculator program (addition, subtraction, multiplication, division).' } else { $notFound += 'pair_7' }

$r = $d.Content
$found = $r.Find.Execute('**Week 3: Strings & String Manipulation**')
if ($found) { $r.Text = '**Week 3: Control Flow (Conditional Statements)**' } else { $notFound += 'pair_8' }

$r = $d.Content
$found = $r.Find.Execute('* **Topic:**  Advanced String Operations.')
if ($found) { $r.Text = '* **Topic:** Controlling program flow using conditional statements.' } else { $notFound += 'pair_9' }

$r = $d.Content
$found = $r.Find.Execute('* **Subtopics:** String methods (upper(), lower(), split(), join(), find(), replace()), string formatting (f-strings), working with escape characters.')
if ($found) { $r.Text = '* **Subtopics:** `if`, `elif`, `else` statements.  Nested conditional statements. Boolean logic.' } else { $notFound += 'pair_10' }

$r = $d.Content
$found = $r.Find.Execute('* **Activities:** Write a program to analyze text (e.g., count word occurrences). Create a program that formats a name and address.')
if ($found) { $r.Text = '* **Activities:**  Write programs with different conditional scenarios (e.g., checking for even/odd numbers, grading system).' } else { $notFound += 'pair_11' }

$r = $d.Content
$found = $r.Find.Execute('**Week 4: Lists, Tuples & Sets**')
if ($found) { $r.Text = '**Week 4: Control Flow (Loops)**' } else { $notFound += 'pair_12' }

$r = $d.Content
$found = $r.Find.Execute('* **Topic:** Collection Data Types')
if ($found) { $r.Text = '* **Topic:** Iterating through data using loops.' } else { $notFound += 'pair_13' }

$r = $d.Content
$found = $r.Find.Execute('* **Subtopics:**  Creating and manipulating lists, tuples, and sets.  List comprehensions (introduction).  Differences between mutable and immutable data types.')
if ($found) { $r.Text = '* **Subtopics:** `for` loops (iterating through lists, ranges). `while` loops.  `break` and `continue` statements.  Nested loops.' } else { $notFound += 'pair_14' }

$r = $d.Content
$found = $r.Find.Execute('* **Activities:**  Write programs that use lists to store and manipulate data.  Compare the performance of lists, tuples and sets in different scenarios.')
if ($found) { $r.Text = '* **Activities:**  Write programs that use loops to perform repetitive tasks (e.g., printing patterns, calculating factorials).' } else { $notFound += 'pair_15' }

$r = $d.Content
$found = $r.Find.Execute('**Week 5: Dictionaries & Control Flow (Part 1)**')
if ($found) { $r.Text = '**Week 5: Data Structures I (Lists & Tuples)**' } else { $notFound += 'pair_16' }

$r = $d.Content
$found = $r.Find.Execute('* **Topic:**  Data Organization & Decision Making.')
if ($found) { $r.Text = '* **Topic:** Working with lists and tuples.' } else { $notFound += 'pair_17' }

$r = $d.Content
$found = $r.Find.Execute('* **Subtopics:**  Creating and accessing dictionaries, iterating through dictionaries.  Conditional statements: `if`, `elif`, `else`.')
if ($found) { $r.Text = '* **Subtopics:**  Creating, accessing, and modifying lists and tuples. List comprehension.  Slicing.  Methods for lists and tuples.' } else { $notFound += 'pair_18' }

$r = $d.Content
$found = $r.Find.Execute('* **Activities:**  Create a program to store and access student information using a dictionary.  Write a program that uses conditional statements to determine grades based on scores.')
if ($found) { $r.Text = '* **Activities:**  Practice manipulating lists and tuples. Create a program to manage a list of student names.' } else { $notFound += 'pair_19' }

$r = $d.Content
$found = $r.Find.Execute('**Week 6: Control Flow (Part 2) & Loops**')
if ($found) { $r.Text = '**Week 6: Data Structures II (Dictionaries & Sets)**' } else { $notFound += 'pair_20' }

$r = $d.Content
$found = $r.Find.Execute('* **Topic:**  Iteration and Control.')
if ($found) { $r.Text = '* **Topic:**  Working with dictionaries and sets.' } else { $notFound += 'pair_21' }

$r = $d.Content
$found = $r.Find.Execute('* **Subtopics:** `for` loops, `while` loops, `break` and `continue` statements, nested loops.')
if ($found) { $r.Text = '* **Subtopics:**  Creating, accessing, and modifying dictionaries.  Set operations (union, intersection, difference).' } else { $notFound += 'pair_22' }

$r = $d.Content
$found = $r.Find.Execute('* **Activities:**  Write programs using `for` and `while` loops to accomplish different tasks (e.g., print numbers, calculate sums).  Practice using `break` and `continue` to control loop execution.')
if ($found) { $r.Text = '* **Activities:**  Create a program to store and retrieve student information using a dictionary.  Implement a program to find common elements in two sets.' } else { $notFound += 'pair_23' }

$r = $d.Content
$found = $r.Find.Execute('**Week 7: Functions & Modularity**')
if ($found) { $r.Text = '**Week 7: Functions**' } else { $notFound += 'pair_24' }

$r = $d.Content
$found = $r.Find.Execute('* **Topic:**  Code Reusability.')
if ($found) { $r.Text = '* **Topic:**  Modularizing code with functions.' } else { $notFound += 'pair_25' }

$r = $d.Content
$found = $r.Find.Execute('* **Subtopics:** Defining functions, function arguments, return values, scope, docstrings.')
if ($found) { $r.Text = '* **Subtopics:**  Defining and calling functions.  Parameters and arguments.  Return values.  Scope and lifetime of variables.' } else { $notFound += 'pair_26' }

$r = $d.Content
$found = $r.Find.Execute('* **Activities:**  Write functions to perform specific tasks (e.g., calculate area, check if a number is prime). Create a program that uses multiple functions.')
if ($found) { $r.Text = '* **Activities:**  Write functions to perform specific tasks (e.g., calculating area, checking for prime numbers).' } else { $notFound += 'pair_27' }

$r = $d.Content
$found = $r.Find.Execute('**Week 8: Modules & Packages**')
if ($found) { $r.Text = '**Week 8: Modules and Packages**' } else { $notFound += 'pair_28' }

$r = $d.Content
$found = $r.Find.Execute('* **Topic:**  Working with External Code.')
if ($found) { $r.Text = '* **Topic:** Utilizing built-in and external modules.' } else { $notFound += 'pair_29' }

$r = $d.Content
$found = $r.Find.Execute('* **Subtopics:** Importing modules (e.g., `math`, `random`, `datetime`), using built-in functions, installing packages using `pip`, introduction to a library (e.g., `requests`).')
if ($found) { $r.Text = '* **Subtopics:**  Importing modules.  Using built-in modules (e.g., `math`, `random`, `datetime`). Installing packages using `pip`.' } else { $notFound += 'pair_30' }

$r = $d.Content
$found = $r.Find.Execute('* **Activities:**  Write a program that uses functions from the `math` module. Install a package and use its functionality in a program.')
if ($found) { $r.Text = '* **Activities:**  Use `math` module functions in a program. Install a simple package (e.g., `requests`) and use its functionality.' } else { $notFound += 'pair_31' }

$r = $d.Content
$found = $r.Find.Execute('* **Topic:**  Robust Code.')
if ($found) { $r.Text = '* **Topic:**  Handling errors gracefully.' } else { $notFound += 'pair_32' }

$r = $d.Content
$found = $r.Find.Execute('* **Subtopics:**  `try`, `except`, `finally` blocks, handling specific exceptions.')
if ($found) { $r.Text = '* **Subtopics:**  `try`, `except`, `finally` blocks.  Common exceptions (e.g., `TypeError`, `ValueError`, `FileNotFoundError`).' } else { $notFound += 'pair_33' }

$r = $d.Content
$found = $r.Find.Execute('* **Activities:**  Write a program that gracefully handles potential errors (e.g., file not found, division by zero).')
if ($found) { $r.Text = '* **Activities:**  Write a program that handles potential errors (e.g., division by zero, file not found).' } else { $notFound += 'pair_34' }

$r = $d.Content
$found = $r.Find.Execute('**Week 10: Introduction to Object-Oriented Programming (OOP)**')
if ($found) { $r.Text = '**Week 10: Object-Oriented Programming (OOP) Basics**' } else { $notFound += 'pair_35' }

$r = $d.Content
$found = $r.Find.Execute('* **Topic:**  OOP Concepts.')
if ($found) { $r.Text = '* **Topic:** Introduction to OOP concepts.' } else { $notFound += 'pair_36' }

$r = $d.Content
$found = $r.Find.Execute('* **Subtopics:** Classes, objects, attributes, methods, constructors (`__init__`).')
if ($found) { $r.Text = '* **Subtopics:** Classes and objects.  Attributes and methods.  Constructors (`__init__`).' } else { $notFound += 'pair_37' }

$r = $d.Content
$found = $r.Find.Execute('* **Activities:**  Create a simple class (e.g., a `Dog` class with attributes like name and breed).')
if ($found) { $r.Text = '* **Activities:** Create a simple class (e.g., a `Dog` class with attributes like name and breed).' } else { $notFound += 'pair_38' }

$r = $d.Content
$found = $r.Find.Execute('**Week 11: OOP (continued) & File Handling**')
if ($found) { $r.Text = '**Week 11: OOP (Inheritance & Encapsulation)**' } else { $notFound += 'pair_39' }

$r = $d.Content
$found = $r.Find.Execute('* **Topic:**  Advanced OOP & File I/O.')
if ($found) { $r.Text = '* **Topic:**  Advanced OOP concepts.' } else { $notFound += 'pair_40' }

$r = $d.Content
$found = $r.Find.Execute('* **Subtopics:** Inheritance, encapsulation, polymorphism (basic concepts). Reading and writing files (text files, CSV files).')
if ($found) { $r.Text = '* **Subtopics:** Inheritance (creating subclasses).  Encapsulation (data hiding).' } else { $notFound += 'pair_41' }

$r = $d.Content
$found = $r.Find.Execute('* **Activities:**  Create a program that uses inheritance to extend a class. Write a program that reads and writes data to a file.')
if ($found) { $r.Text = '* **Activities:**  Extend the `Dog` class to create subclasses (e.g., `Labrador`, `GoldenRetriever`).' } else { $notFound += 'pair_42' }

$r = $d.Content
$found = $r.Find.Execute('**Week 12: Mini-Project & Review**')
if ($found) { $r.Text = '**Week 12: File I/O and Mini-Project**' } else { $notFound += 'pair_43' }

$r = $d.Content
$found = $r.Find.Execute('* **Topic:**  Putting it all together.')
if ($found) { $r.Text = '* **Topic:**  Working with files and a culminating project.' } else { $notFound += 'pair_44' }

$r = $d.Content
$found = $r.Find.Execute('* **Subtopics:**  Project brainstorming and implementation (calculator, quiz app, simple data parser, etc.). Review of key concepts.')
if ($found) { $r.Text = '* **Subtopics:**  Reading and writing files.  Working with CSV or JSON data.  Mini-project presentation.' } else { $notFound += 'pair_45' }

$vt = [char]11
$finalOld = '* **Activities:**  Complete a mini-project that incorporates the concepts learned throughout the course.  Prepare for a final assessment (optional).' + $vt + $vt + $vt + 'This lesson plan provides a flexible framework. Adjust the pace and depth of coverage based on the students'' progress and understanding.  Remember to incorporate regular quizzes and coding exercises to reinforce learning.'
$finalNew = '* **Activities:**  Complete a mini-project (e.g., a simple calculator, a quiz app, or a data parser).  Present the project to the class.'
$r = $d.Content
$found = $r.Find.Execute($finalOld)
if ($found) { $r.Text = $finalNew } else { $notFound += 'final_pair' }

Write-Output "notFoundCount=$($notFound.Count)"
if ($notFound.Count -gt 0) { Write-Output ($notFound -join ",") }
Write-Output "done"